# Applies the "Generate Report for Handoff" update to the localization-status
# workbook:
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - Latest Handoff/HO timestamps are bumped a couple of minutes later
#   - Priority changes from "ht" to "mt"
#   - A new Error Detail message is recorded for the 1e4daa63... item
#   - A few report columns are narrowed / widened

$wb = $excel.ActiveWorkbook

$statusText = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3611e5106f5391dafbbef73fea35cf0eab388ac/e2e/1e4daa63-87a3-4c3b-be57-3f2ff07dceb2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a862bab4dbbf842c6116109e5ffedd551b21d947/e2e/1e4daa63-87a3-4c3b-be57-3f2ff07dceb2.md."

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Range("G2").Value = "2016-10-14 08:46:50"
$wsOverview.Range("G3").Value = "2016-10-14 08:46:50"

$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

$wsZhCn.Range("H2").Value = "2016-10-14 08:46:39"
$wsZhCn.Range("H3").Value = "2016-10-14 08:46:39"

$wsZhCn.Range("P2").Value = $errorDetail

$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

$wsDeDe.Range("H2").Value = "2016-10-14 08:46:50"
$wsDeDe.Range("H3").Value = "2016-10-14 08:46:50"

$wsDeDe.Range("P2").Value = $errorDetail

$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
